$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 395
$ws.Range("F5").Value = 188
$ws.Range("F6").Value = 13
$ws.Range("F7").Value = 774
$ws.Range("F8").Value = 74
$ws.Range("F9").Value = 9933
$ws.Range("F11").Value = 2773
$ws.Range("F13").Value = 2412
$ws.Range("F14").Value = 2716
$ws.Range("F16").Value = 295
$ws.Range("F17").Value = 2115
$ws.Range("F19").Value = 88
$ws.Range("F22").Value = 98
$ws.Range("F25").Value = 182
$ws.Range("F26").Value = 607
$ws.Range("F27").Value = 1302
$ws.Range("F29").Value = 97
$ws.Range("F32").Value = 1852
$ws.Range("F33").Value = 2883
$ws.Range("F36").Value = 1012
$ws.Range("F37").Value = 367
$ws.Range("F39").Value = 1282
$ws.Range("F40").Value = 69
$ws.Range("F41").Value = 91
$ws.Range("F42").Value = 60
$ws.Range("F43").Value = 28
$ws.Range("F44").Value = 33

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 165
$ws.Range("F6").Value = 15
$ws.Range("F7").Value = 2
$ws.Range("F15").Value = 165

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 964
$ws.Range("F5").Value = 1916

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 964
$ws.Range("F6").Value = 395
$ws.Range("F9").Value = 188
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 774
$ws.Range("F12").Value = 74
$ws.Range("F13").Value = 9933
$ws.Range("F16").Value = 2774
$ws.Range("F18").Value = 2412
$ws.Range("F19").Value = 2716
$ws.Range("F20").Value = 295
$ws.Range("F21").Value = 2115
$ws.Range("F23").Value = 88
$ws.Range("F28").Value = 182
$ws.Range("F29").Value = 607
$ws.Range("F30").Value = 1302
$ws.Range("F34").Value = 1852
$ws.Range("F36").Value = 2883
$ws.Range("F37").Value = 1012
$ws.Range("F39").Value = 367
$ws.Range("F44").Value = 1282
$ws.Range("F45").Value = 69
$ws.Range("F46").Value = 60
$ws.Range("F47").Value = 28
$ws.Range("F48").Value = 33
$ws.Range("F49").Value = 165
